$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9437782764434814
$ws.Range("B1").Value = 1.953955292701721
$ws.Range("C1").Value = 4.339358329772949
$ws.Range("D1").Value = 3.269311666488647
$ws.Range("E1").Value = 1.442743301391602
